$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, shifting existing rows 141:215 down to 142:216.
$ws.Range("A141").EntireRow.Insert()

# Populate the newly inserted row 141 with the new weekly data point.
$ws.Range("A141").Value = 5
$ws.Range("B141").Value = "Macroferia Regional de Talca"
$ws.Range("C141").Value = "Maule"
$ws.Range("D141").Value = 45016
$ws.Range("E141").Value = 7
$ws.Range("F141").Value = 100112030
$ws.Range("G141").Value = "Poroto granado"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 200
$ws.Range("K141").Value = 35000
$ws.Range("L141").Value = 35000
$ws.Range("M141").Value = 35000
$ws.Range("N141").Value = "`$/saco 25 kilos"
$ws.Range("O141").Value = "Región del Maule"
$ws.Range("P141").Value = 1400
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"
